$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing numeric data cells (rows 2-11) with the latest nowcast run ---
foreach ($pair in @(
    @{Addr="B2"; Val=0.24080207187219643}
    @{Addr="C2"; Val=0}
    @{Addr="D2"; Val=0}
    @{Addr="E2"; Val=0}
    @{Addr="F2"; Val=0}
    @{Addr="G2"; Val=0}
    @{Addr="H2"; Val=0}
    @{Addr="I2"; Val=0}
    @{Addr="J2"; Val=0}
    @{Addr="K2"; Val=0}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B3"; Val=0.28436211342840439}
    @{Addr="C3"; Val=0}
    @{Addr="D3"; Val=0.0012439662374659714}
    @{Addr="E3"; Val=-0.0037817599196970104}
    @{Addr="F3"; Val=-0.0025283334682744876}
    @{Addr="G3"; Val=0.015146769417247196}
    @{Addr="H3"; Val=-0.00057838198981614408}
    @{Addr="I3"; Val=0.033991010788495674}
    @{Addr="J3"; Val=0}
    @{Addr="K3"; Val=0.000066770490786766112}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B4"; Val=0.28208020537492945}
    @{Addr="C4"; Val=-0.0035241831576743455}
    @{Addr="D4"; Val=0}
    @{Addr="E4"; Val=-0.0015560534800836596}
    @{Addr="F4"; Val=-0.000078031851677072832}
    @{Addr="G4"; Val=0}
    @{Addr="H4"; Val=-0.0019327886606760635}
    @{Addr="I4"; Val=0.0045734133020846612}
    @{Addr="J4"; Val=-0.000019697362304567411}
    @{Addr="K4"; Val=0.00025543315685611079}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B5"; Val=0.2118210139814253}
    @{Addr="C5"; Val=0.00017098688669906457}
    @{Addr="D5"; Val=-0.021669175339521635}
    @{Addr="E5"; Val=-0.0061503162577170624}
    @{Addr="F5"; Val=0.0031600996313077712}
    @{Addr="G5"; Val=-0.055992514771728279}
    @{Addr="H5"; Val=-0.00062809303110503947}
    @{Addr="I5"; Val=0.010736702060784356}
    @{Addr="J5"; Val=0}
    @{Addr="K5"; Val=0.00011311942777669781}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B6"; Val=0.37974679247291698}
    @{Addr="C6"; Val=0.21468273181312142}
    @{Addr="D6"; Val=0}
    @{Addr="E6"; Val=0.0033733487000692638}
    @{Addr="F6"; Val=0.00018191970818550758}
    @{Addr="G6"; Val=0}
    @{Addr="H6"; Val=0.00070713972826768798}
    @{Addr="I6"; Val=-0.051037856639376288}
    @{Addr="J6"; Val=0}
    @{Addr="K6"; Val=0.000018495181224098545}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B7"; Val=0.38619987317568699}
    @{Addr="C7"; Val=0}
    @{Addr="D7"; Val=-0.040400119627963769}
    @{Addr="E7"; Val=0.011964863103795773}
    @{Addr="F7"; Val=0.012449969319882239}
    @{Addr="G7"; Val=0.015643289302739314}
    @{Addr="H7"; Val=0}
    @{Addr="I7"; Val=0.0047611922848377463}
    @{Addr="J7"; Val=0}
    @{Addr="K7"; Val=0.002033886319478706}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B8"; Val=0.31795258949313859}
    @{Addr="C8"; Val=0.0037618529577373543}
    @{Addr="D8"; Val=0}
    @{Addr="E8"; Val=-0.015389088370375149}
    @{Addr="F8"; Val=-0.015394683700806744}
    @{Addr="G8"; Val=0}
    @{Addr="H8"; Val=-0.000032896263434469691}
    @{Addr="I8"; Val=-0.04307682182660006}
    @{Addr="J8"; Val=0}
    @{Addr="K8"; Val=0.0018843535209306927}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B9"; Val=-0.071497957548829705}
    @{Addr="C9"; Val=0}
    @{Addr="D9"; Val=-0.23930056372694083}
    @{Addr="E9"; Val=-0.0064611883494527936}
    @{Addr="F9"; Val=-0.13977766188943633}
    @{Addr="G9"; Val=-0.0087111813587890073}
    @{Addr="H9"; Val=-0.0017849062377760128}
    @{Addr="I9"; Val=0.0047681202030623376}
    @{Addr="J9"; Val=0}
    @{Addr="K9"; Val=0.0018168343173642909}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B10"; Val=0.054566760516671781}
    @{Addr="C10"; Val=0.17759045202830737}
    @{Addr="D10"; Val=0}
    @{Addr="E10"; Val=0.0020123662740876092}
    @{Addr="F10"; Val=0.0018983461025653334}
    @{Addr="G10"; Val=0}
    @{Addr="H10"; Val=-0.00053475956616136713}
    @{Addr="I10"; Val=0.040294149890125641}
    @{Addr="J10"; Val=-0.10488631362318888}
    @{Addr="K10"; Val=0.0096904769597657703}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

foreach ($pair in @(
    @{Addr="B11"; Val=0.30558649411408689}
    @{Addr="C11"; Val=0}
    @{Addr="D11"; Val=0.41040944622830727}
    @{Addr="E11"; Val=-0.08157479969732763}
    @{Addr="F11"; Val=-0.17758697045414956}
    @{Addr="G11"; Val=-0.033618800681720573}
    @{Addr="H11"; Val=-0.00078157726851500135}
    @{Addr="I11"; Val=0.10183556307617514}
    @{Addr="J11"; Val=0}
    @{Addr="K11"; Val=0.032336872394645472}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

# --- Append new row 12 (2025-08-30) ---
$a12 = $ws.Range("A12")
$a12.Value = "'2025-08-30"
$a12.Style = "Normal"

foreach ($pair in @(
    @{Addr="B12"; Val=0.086829374710636814}
    @{Addr="C12"; Val=-0.12507595826134951}
    @{Addr="D12"; Val=0}
    @{Addr="E12"; Val=-0.0019631484652427259}
    @{Addr="F12"; Val=-0.0018500681955330132}
    @{Addr="G12"; Val=0}
    @{Addr="H12"; Val=-0.00024982997216308985}
    @{Addr="I12"; Val=-0.042613970191810242}
    @{Addr="J12"; Val=0}
    @{Addr="K12"; Val=-0.047004144317351504}
)) { $ws.Range($pair.Addr).Value = $pair.Val }

# --- Column width tweaks (engine snaps ColumnWidth to sixths of a character; closest achievable) ---
$ws.Columns.Item(3).ColumnWidth = 14.833333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
